$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''30.609.13'
$ws.Range('E2').Value = '''  +1.69%  '
$ws.Range('D3').Value = '''1.923.49'
$ws.Range('E3').Value = '''  +4.17%  '
$ws.Range('D4').Value = '''1.0000'
$ws.Range('E4').Value = '''  -0.19%  '
$ws.Range('D5').Value = '''248.09'
$ws.Range('E5').Value = '''  +5.62%  '
$ws.Range('E6').Value = '''  -0.08%  '
$ws.Range('D7').Value = '''0.4736'
$ws.Range('E7').Value = '''  +2.04%  '
$ws.Range('D8').Value = '''0.2913'
$ws.Range('E8').Value = '''  +4.73%  '
$ws.Range('E9').Value = '''  +6.61%  '
$ws.Range('D10').Value = '''105.68'
$ws.Range('E10').Value = '''  +9.52%  '
$ws.Range('D11').Value = '''18.44'
$ws.Range('E11').Value = '''  +1.23%  '
$ws.Range('B12').Value = 'TRON'
$ws.Range('C12').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D12').Value = '''0.07720'
$ws.Range('E12').Value = '''  +2.44%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = '''1.907.00'
$ws.Range('E13').Value = '''  +3.64%  '
$ws.Range('D14').Value = '''5.303'
$ws.Range('E14').Value = '''  +7.00%  '
$ws.Range('D15').Value = '''0.6734'
$ws.Range('E15').Value = '''  +7.50%  '
$ws.Range('D16').Value = '''286.76'
$ws.Range('E16').Value = '''  -2.54%  '
$ws.Range('D17').Value = '''30.618.90'
$ws.Range('E17').Value = '''  +1.72%  '
$ws.Range('D18').Value = '''0.000007647'
$ws.Range('E18').Value = '''  +4.12%  '
$ws.Range('B19').Value = 'Avalanche'
$ws.Range('C19').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D19').Value = '''12.96'
$ws.Range('E19').Value = '''  +2.58%  '
$ws.Range('B20').Value = 'Dai'
$ws.Range('C20').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D20').Value = '''0.9996'
$ws.Range('E20').Value = '''  -0.32%  '
$ws.Range('D21').Value = '''2.161.98'
$ws.Range('E21').Value = '''  +3.86%  '
$ws.Range('D22').Value = '''5.425'
$ws.Range('E22').Value = '''  +8.73%  '
$ws.Range('E23').Value = '''  -0.35%  '
$ws.Range('D24').Value = '''6.310'
$ws.Range('E24').Value = '''  +4.47%  '
$ws.Range('D25').Value = '''9.426'
$ws.Range('E25').Value = '''  +3.80%  '
$ws.Range('D26').Value = '''168.40'
$ws.Range('E26').Value = '''  +2.24%  '
$ws.Range('D27').Value = '''20.87'
$ws.Range('E27').Value = '''  +8.57%  '
$ws.Range('D28').Value = '''2.137'
$ws.Range('E28').Value = '''  +11.59%  '
$ws.Range('D29').Value = '''0.1082'
$ws.Range('E29').Value = '''  +0.55%  '
$ws.Range('D30').Value = '''1.370'
$ws.Range('E30').Value = '''  +3.11%  '
$ws.Range('D31').Value = '''4.196'
$ws.Range('E31').Value = '''  +4.86%  '
$ws.Range('D32').Value = '''4.145'
$ws.Range('E32').Value = '''  +9.04%  '
$ws.Range('D33').Value = '''0.05046'
$ws.Range('E33').Value = '''  +3.01%  '
$ws.Range('D34').Value = '''0.7429'
$ws.Range('E34').Value = '''  +2.67%  '
$ws.Range('D35').Value = '''1.159'
$ws.Range('E35').Value = '''  +4.86%  '
$ws.Range('D36').Value = '''0.02083'
$ws.Range('E36').Value = '''  +8.98%  '
$ws.Range('D37').Value = '''2.745'
$ws.Range('E37').Value = '''  +0.45%  '
$ws.Range('D38').Value = '''2.692'
$ws.Range('E38').Value = '''  +1.22%  '
$ws.Range('D39').Value = '''2.069'
$ws.Range('E39').Value = '''  +5.54%  '
$ws.Range('D40').Value = '''111.11'
$ws.Range('E40').Value = '''  +6.01%  '
$ws.Range('D41').Value = '''0.8844'
$ws.Range('E41').Value = '''  +2.59%  '
$ws.Range('D42').Value = '''0.4390'
$ws.Range('E42').Value = '''  +8.95%  '
$ws.Range('D43').Value = '''5.946'
$ws.Range('E43').Value = '''  +5.82%  '
$ws.Range('D44').Value = '''0.9998'
$ws.Range('E44').Value = '''  -0.07%  '
$ws.Range('D45').Value = '''67.42'
$ws.Range('E45').Value = '''  +3.82%  '
$ws.Range('D46').Value = '''7.268'
$ws.Range('E46').Value = '''  +3.79%  '
$ws.Range('D47').Value = '''9.330'
$ws.Range('E47').Value = '''  +4.33%  '
$ws.Range('D48').Value = '''47.96'
$ws.Range('E48').Value = '''  +17.73%  '
$ws.Range('D49').Value = '''0.1238'
$ws.Range('E49').Value = '''  +4.53%  '
$ws.Range('D50').Value = '''35.08'
$ws.Range('E50').Value = '''  +3.36%  '
$ws.Range('D51').Value = '''0.4078'
$ws.Range('E51').Value = '''  +10.09%  '
